$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 125; this shifts existing rows 125-180 down to 126-181
$ws.Rows("125").Insert()

# Populate the newly inserted row 125 with the new data record
$ws.Range("A125").Value = 9
$ws.Range("B125").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = 44510
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 300000001
$ws.Range("G125").Value = "Rabanito"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 7900
$ws.Range("K125").Value = 2500
$ws.Range("L125").Value = 3000
$ws.Range("M125").Value = 2747
$ws.Range("N125").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O125").Value = "Provincia de Chacabuco"
$ws.Range("P125").Value = 27
$ws.Range("Q125").Value = 100
$ws.Range("R125").Value = "Hortaliza"
